$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# The document currently ends with two trailing paragraphs:
#   <w:p/>                         (empty paragraph)
#   <w:p><w:r><w:tab/></w:r></w:p> (paragraph containing just a tab)
# These are replaced by a short "Things that need attention:" paragraph
# followed by two new bulleted ("ListParagraph"/numId 1) list items.

# 1) Turn the empty paragraph into "Things that" + " need attention:" (two runs)
$p = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$xml = "<w:p $wNs><w:r><w:t>Things that</w:t></w:r><w:r><w:t xml:space=`"preserve`"> need attention:</w:t></w:r></w:p>"
[void]$p.Range.InsertXML($xml)

# 2) Turn the tab-only paragraph into the first bullet item: "Scene" + ".cpp (missing two methods)"
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$xml = "<w:p $wNs><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr></w:pPr><w:r><w:t>Scene</w:t></w:r><w:r><w:t>.cpp (missing two methods)</w:t></w:r></w:p>"
[void]$p.Range.InsertXML($xml)

# 3) Append a new bullet item after it: "Memory management (how to handle it?)"
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
[void]$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$xml = "<w:p $wNs><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr></w:pPr><w:r><w:t>Memory management (how to handle it?)</w:t></w:r></w:p>"
[void]$p.Range.InsertXML($xml)
